$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) & 2) Swap the w:name values of the two duplicate-bookmark paragraphs.
#    Each of these paragraphs is otherwise empty (no visible text), so we
#    rebuild its XML in place via InsertXML (the replacement Range exactly
#    spans the existing paragraph, so no paragraph split/merge happens and
#    neighbouring paragraphs are left untouched). Because InsertXML silently
#    drops paragraph-level formatting that happens to equal the resolved
#    style defaults, we re-assert the (redundant but originally explicit)
#    alignment / line-spacing afterwards through the Format object so the
#    paragraph mark keeps its direct <w:spacing>/<w:jc>.
# ---------------------------------------------------------------------------

function Swap-BookmarkParagraph($paraIndex, $id1, $name1, $id2, $name2, $id3, $name3, $id4, $name4) {
    $p = $d.Paragraphs.Item($paraIndex)
    $r = $p.Range

    $frag = '<w:p><w:pPr><w:pStyle w:val="Normal"/><w:spacing w:lineRule="auto" w:line="360"/>' + `
        '<w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/>' + `
        '<w:color w:val="000000"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' + `
        '<w:r><w:rPr><w:rFonts w:cs="Arial" w:ascii="Arial" w:hAnsi="Arial"/><w:color w:val="000000"/>' + `
        '<w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:r>' + `
        '<w:bookmarkStart w:id="' + $id1 + '" w:name="' + $name1 + '"/>' + `
        '<w:bookmarkStart w:id="' + $id2 + '" w:name="' + $name2 + '"/>' + `
        '<w:bookmarkStart w:id="' + $id3 + '" w:name="' + $name3 + '"/>' + `
        '<w:bookmarkStart w:id="' + $id4 + '" w:name="' + $name4 + '"/>' + `
        '<w:bookmarkEnd w:id="' + $id3 + '"/><w:bookmarkEnd w:id="' + $id4 + '"/></w:p>'

    $newXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body>' + $frag + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $r.InsertXML($newXml)

    $p2 = $d.Paragraphs.Item($paraIndex)
    $p2.Format.Alignment = 3        # wdAlignParagraphJustify  -> <w:jc w:val="both"/>
    $p2.Format.LineSpacingRule = 1  # wdLineSpace1pt5          -> <w:spacing w:line="360" w:lineRule="auto"/>
}

# Paragraph containing ids 1-4 / _Hlk157071468 & _Hlk162862782 (swap names).
Swap-BookmarkParagraph 17 1 "_Hlk162862782" 2 "_Hlk157071468" 3 "_Hlk162862782" 4 "_Hlk157071468"

# Paragraph containing ids 5-8 / _Hlk1570714681 & _Hlk1628627821 (swap names).
Swap-BookmarkParagraph 21 5 "_Hlk1628627821" 6 "_Hlk1570714681" 7 "_Hlk1628627821" 8 "_Hlk1570714681"

# ---------------------------------------------------------------------------
# 3) Split the {rua} run (the one following "... situado na ") into
#    "{" + "referencia1" + "}", giving the middle piece an explicit black
#    color while leaving the braces in the original (colorless) formatting.
# ---------------------------------------------------------------------------

$full = $d.Content
$full.Find.Execute("situado na {rua};", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$ruaStart = $full.Start + "situado na ".Length
$ruaEnd = $ruaStart + "{rua}".Length

$rMid = $d.Range($ruaStart + 1, $ruaEnd - 1)   # just the "rua" text
$rMid.Text = "referencia1"

$midLen = "referencia1".Length
$rMid2 = $d.Range($ruaStart + 1, $ruaStart + 1 + $midLen)
$rMid2.Font.Color = 0   # wdColorBlack -> <w:color w:val="000000"/>
